$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 45
$ws.Range("I58").Value = 45
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 135
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 15
$ws.Range("N58").ClearContents()

$ws.Range("H62").Value = 8409466
$ws.Range("I62").Value = 14293912
$ws.Range("J62").Value = 3115.5715
$ws.Range("K62").Value = 14293912
$ws.Range("L62").Value = 3115.5715
$ws.Range("M62").Value = -14293288
$ws.Range("N62").Value = -4363.5715

$ws.Range("H64").Value = 111114450
$ws.Range("I64").Value = 111114450
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 111114450
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -111114202

$ws.Range("H65").Value = 8409466
$ws.Range("I65").Value = 14293912
$ws.Range("J65").Value = 3115.5715
$ws.Range("K65").Value = 71469560
$ws.Range("L65").Value = 15577.8575
$ws.Range("M65").Value = -71466440
$ws.Range("N65").Value = -21817.8575

$ws.Range("H67").Value = 111114450
$ws.Range("I67").Value = 111114450
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 111114450
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -111113592

$ws.Range("H103").Value = 282.3684
$ws.Range("I103").Value = 268.13333
$ws.Range("J103").Value = 335.75
$ws.Range("K103").Value = 804.39999
$ws.Range("L103").Value = 1007.25
$ws.Range("M103").Value = -218.39999

$ws.Range("H107").Value = 2717.2632
$ws.Range("I107").Value = 1487.6923
$ws.Range("J107").Value = 5381.3335
$ws.Range("K107").Value = 1487.6923
$ws.Range("L107").Value = 5381.3335
$ws.Range("M107").Value = 432.3077000000001
$ws.Range("N107").Value = -9221.333500000001

$ws.Range("H113").Value = 25271.715
$ws.Range("I113").Value = 5001.3335
$ws.Range("J113").Value = 40474.5
$ws.Range("K113").Value = 5001.3335
$ws.Range("L113").Value = 40474.5
$ws.Range("M113").Value = -1747.3335

$ws.Range("H132").Value = 498541.12
$ws.Range("I132").Value = 646058.5
$ws.Range("J132").Value = 6816.5
$ws.Range("K132").Value = 1938175.5
$ws.Range("L132").Value = 20449.5
$ws.Range("M132").Value = -1935645.5

$ws.Range("H138").Value = 3161.5103
$ws.Range("I138").Value = 760.8
$ws.Range("J138").Value = 5662.25
$ws.Range("K138").Value = 2282.4
$ws.Range("L138").Value = 16986.75
$ws.Range("M138").Value = 2857.6
$ws.Range("N138").Value = -27266.75

$ws.Range("H141").Value = 2802.0833
$ws.Range("I141").Value = 2836.9565
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 8510.869499999999
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -3330.869499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7447.5386
$ws.Range("I2").Value = 4082.2
$ws.Range("J2").Value = 18665.334
$ws.Range("K2").Value = 4082.2
$ws.Range("L2").Value = 18665.334
$ws.Range("M2").Value = -3969.2

$ws.Range("H32").Value = 2780526
$ws.Range("I32").Value = 1122.3214
$ws.Range("J32").Value = 9265801
$ws.Range("K32").Value = 1122.3214
$ws.Range("L32").Value = 9265801
$ws.Range("M32").Value = -835.3214

$ws.Range("H45").Value = 3018.6667
$ws.Range("I45").Value = 2450.2222
$ws.Range("J45").Value = 3587.111
$ws.Range("K45").Value = 2450.2222
$ws.Range("L45").Value = 3587.111
$ws.Range("M45").Value = -2073.2222
$ws.Range("N45").Value = -4341.111

$ws.Range("H110").Value = 4845.125
$ws.Range("I110").Value = 3361.5293
$ws.Range("J110").Value = 6526.533
$ws.Range("K110").Value = 3361.5293
$ws.Range("L110").Value = 6526.533
$ws.Range("M110").Value = -1316.5293

$ws.Range("H116").Value = 7447.5386
$ws.Range("I116").Value = 4082.2
$ws.Range("J116").Value = 18665.334
$ws.Range("K116").Value = 4082.2
$ws.Range("L116").Value = 18665.334
$ws.Range("M116").Value = -1788.2

$ws.Range("H132").Value = 2085213.4
$ws.Range("I132").Value = 4330661
$ws.Range("J132").Value = 120446.75
$ws.Range("K132").Value = 12991983
$ws.Range("L132").Value = 361340.25
$ws.Range("M132").Value = -12989453
$ws.Range("N132").Value = -366400.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7447.5386
$ws.Range("I3").Value = 4082.2
$ws.Range("J3").Value = 18665.334
$ws.Range("K3").Value = 4082.2
$ws.Range("L3").Value = 18665.334
$ws.Range("M3").Value = -3968.2

$ws.Range("H10").Value = 463.33334
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 463.33334
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 463.33334
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -743.33334

$ws.Range("H35").Value = 52411
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 52411
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 52411
$ws.Range("N35").Value = -53031

$ws.Range("H99").Value = 5690.4
$ws.Range("I99").Value = 2684.1333
$ws.Range("J99").Value = 8696.666999999999
$ws.Range("K99").Value = 2684.1333
$ws.Range("L99").Value = 8696.666999999999
$ws.Range("M99").Value = -1186.1333
$ws.Range("N99").Value = -11692.667

$ws.Range("H105").Value = 3708.4443
$ws.Range("I105").Value = 7939.6
$ws.Range("J105").Value = 2746.818
$ws.Range("K105").Value = 7939.6
$ws.Range("L105").Value = 2746.818
$ws.Range("M105").Value = -6192.6
$ws.Range("N105").Value = -6240.818

$ws.Range("H107").Value = 20002798
$ws.Range("I107").Value = 25002248
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 25002248
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -25000328

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1677.7
$ws.Range("I107").Value = 1038.2
$ws.Range("J107").Value = 2317.2
$ws.Range("K107").Value = 1038.2
$ws.Range("L107").Value = 2317.2
$ws.Range("M107").Value = 881.8
$ws.Range("N107").Value = -6157.2

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H134").Value = 33343548
$ws.Range("I134").Value = 142863650
$ws.Range("J134").Value = 11345.739
$ws.Range("K134").Value = 428590950
$ws.Range("L134").Value = 34037.217
$ws.Range("M134").Value = -428588415

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1402
$ws.Range("I44").Value = 1402
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 4206
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -3808

$ws.Range("H131").Value = 37682932
$ws.Range("I131").Value = 53333788
$ws.Range("J131").Value = 25643814
$ws.Range("K131").Value = 160001364
$ws.Range("L131").Value = 76931442
$ws.Range("M131").Value = -159996324
$ws.Range("N131").Value = -76941522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1657.1666
$ws.Range("I107").Value = 989.6
$ws.Range("J107").Value = 4995
$ws.Range("K107").Value = 989.6
$ws.Range("L107").Value = 4995
$ws.Range("M107").Value = 930.4
$ws.Range("N107").Value = -8835

$ws.Range("H113").Value = 6805.0625
$ws.Range("I113").Value = 3058.2
$ws.Range("J113").Value = 8508.182000000001
$ws.Range("K113").Value = 3058.2
$ws.Range("L113").Value = 8508.182000000001
$ws.Range("M113").Value = -888.1999999999998
$ws.Range("N113").Value = -12848.182

$ws.Range("H122").Value = 7311.3335
$ws.Range("I122").Value = 4554.1113
$ws.Range("J122").Value = 15583
$ws.Range("K122").Value = 13662.3339
$ws.Range("L122").Value = 46749
$ws.Range("M122").Value = -11212.3339
$ws.Range("N122").Value = -51649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 814.5
$ws.Range("I22").Value = 524.6667
$ws.Range("J22").Value = 1249.25
$ws.Range("K22").Value = 524.6667
$ws.Range("L22").Value = 1249.25
$ws.Range("M22").Value = -229.6667

$ws.Range("H27").Value = 814.5
$ws.Range("I27").Value = 524.6667
$ws.Range("J27").Value = 1249.25
$ws.Range("K27").Value = 524.6667
$ws.Range("L27").Value = 1249.25
$ws.Range("M27").Value = -417.6667

$ws.Range("H46").Value = 19232754
$ws.Range("I46").Value = 786.5625
$ws.Range("J46").Value = 50003900
$ws.Range("K46").Value = 786.5625
$ws.Range("L46").Value = 50003900
$ws.Range("M46").Value = -598.5625
$ws.Range("N46").Value = -50004276

$ws.Range("H61").Value = 3776.6667
$ws.Range("I61").Value = 1821.6957
$ws.Range("J61").Value = 8273.1
$ws.Range("K61").Value = 1821.6957
$ws.Range("L61").Value = 8273.1
$ws.Range("M61").Value = -1619.6957
$ws.Range("N61").Value = -8677.1

$ws.Range("H100").Value = 1982.25
$ws.Range("I100").Value = 1752.6
$ws.Range("J100").Value = 2042.6842
$ws.Range("K100").Value = 1752.6
$ws.Range("L100").Value = 2042.6842
$ws.Range("M100").Value = -1211.6
$ws.Range("N100").Value = -3124.6842

$ws.Range("H113").Value = 3776.6667
$ws.Range("I113").Value = 1821.6957
$ws.Range("J113").Value = 8273.1
$ws.Range("K113").Value = 1821.6957
$ws.Range("L113").Value = 8273.1
$ws.Range("M113").Value = 348.3043
$ws.Range("N113").Value = -12613.1

$ws.Range("H122").Value = 4894.5
$ws.Range("I122").Value = 3225.6924
$ws.Range("J122").Value = 6036.316
$ws.Range("K122").Value = 9677.0772
$ws.Range("L122").Value = 18108.948
$ws.Range("M122").Value = -7227.0772
$ws.Range("N122").Value = -23008.948

$ws.Range("H136").Value = 12825265
$ws.Range("I136").Value = 100000856
$ws.Range("J136").Value = 5324.5
$ws.Range("K136").Value = 300002568
$ws.Range("L136").Value = 15973.5
$ws.Range("M136").Value = -300000018
$ws.Range("N136").Value = -21073.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 5004997.5
$ws.Range("I11").Value = 10000000
$ws.Range("J11").Value = 9995
$ws.Range("K11").Value = 10000000
$ws.Range("L11").Value = 9995
$ws.Range("M11").Value = -9999858

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H113").Value = 23810708
$ws.Range("I113").Value = 51589216
$ws.Range("J113").Value = 558.1429000000001
$ws.Range("K113").Value = 154767648
$ws.Range("L113").Value = 1674.4287
$ws.Range("M113").Value = -154765478

$ws.Range("H122").Value = 2905.56
$ws.Range("I122").Value = 2489.6428
$ws.Range("J122").Value = 3434.9092
$ws.Range("K122").Value = 7468.928400000001
$ws.Range("L122").Value = 10304.7276
$ws.Range("M122").Value = -5018.928400000001
